# Needle calibration data was re-sorted in ascending chronological order (column A = time).
# Rows 5 and 8 already happened to be in the correct position, so we only rewrite the rows
# whose contents actually change, leaving untouched cells with their original serialization.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 8
$numRows = $lastDataRow - $firstDataRow + 1
$numCols = 4

# Read the existing data rows into memory first (row-major: index 0 => sheet row 2, etc.)
$data = New-Object 'object[,]' $numRows, $numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $data[$r, $c] = $ws.Cells.Item($r + $firstDataRow, $c + 1).Value2
    }
}

# Compute the row order (indices into $data) sorted ascending by column A (time),
# using a manual selection sort since Sort-Object's custom expressions are unreliable here.
$order = New-Object 'int[]' $numRows
for ($i = 0; $i -lt $numRows; $i++) { $order[$i] = $i }

for ($i = 0; $i -lt $numRows - 1; $i++) {
    $minIdx = $i
    for ($j = $i + 1; $j -lt $numRows; $j++) {
        if ([double]$data[$order[$j], 0] -lt [double]$data[$order[$minIdx], 0]) {
            $minIdx = $j
        }
    }
    if ($minIdx -ne $i) {
        $tmp = $order[$i]
        $order[$i] = $order[$minIdx]
        $order[$minIdx] = $tmp
    }
}

# Only rewrite rows whose source position actually differs from the destination position,
# so rows that are already correctly placed keep their original on-disk formatting.
for ($i = 0; $i -lt $numRows; $i++) {
    $srcIdx = $order[$i]
    if ($srcIdx -ne $i) {
        $destRow = $i + $firstDataRow
        for ($c = 0; $c -lt $numCols; $c++) {
            $ws.Cells.Item($destRow, $c + 1).Value = $data[$srcIdx, $c]
        }
    }
}

$wb.Save()
